$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E (shifts old E..I to F..J)
$ws.Columns.Item(5).Insert()

# Update D1 text (ID Kategori -> ID Sub Kategori) and set new E1 header
$ws.Range("D1").Value = "ID Sub Kategori"
$ws.Range("E1").Value = "Nama Sub Kategori"

# Remove the old formatted blank row 2
$ws.Rows.Item(2).Delete()

# Add new trailing headers (order chosen so shared-string table indices line up)
$ws.Range("L1").Value = "Garansi IMEI"
$ws.Range("M1").Value = "PPN 11%"
$ws.Range("K1").Value = "Garansi Produk"

# Match header style (bold) like the rest of row 1
$ws.Range("K1:M1").Font.Bold = $true

# Reflect the cursor position left behind in the source workbook
$ws.Range("G7").Select()
